$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.577136993408203
$ws.Range("B1").Value = 1.824631929397583
$ws.Range("C1").Value = 1.892598152160645
$ws.Range("D1").Value = 2.329355716705322
$ws.Range("E1").Value = 3.412607669830322
